$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.074.36'
$ws.Range('E2').Value = '  -0.11%  '

$ws.Range('D3').Value = '3.366.11'
$ws.Range('E3').Value = '  +2.43%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = "'571.18"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '

$ws.Range('D6').Value = "'135.52"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.72%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').Value = '3.366.07'
$ws.Range('E8').Value = '  +2.54%  '

$ws.Range('E9').Value = '  +1.20%  '

$ws.Range('E10').Value = '  +6.43%  '

$ws.Range('D11').Value = "'0.123"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.64%  '

$ws.Range('D12').Value = "'0.390"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.63%  '

$ws.Range('D13').Value = '3.938.41'
$ws.Range('E13').Value = '  +1.86%  '

$ws.Range('E14').Value = '  +2.80%  '

$ws.Range('D15').Value = "'0.0000172"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.34%  '

$ws.Range('D16').Value = '3.350.90'
$ws.Range('E16').Value = '  +1.46%  '

$ws.Range('D17').Value = "'25.12"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.38%  '

$ws.Range('D18').Value = '61.118.33'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('E19').Value = '  +8.40%  '

$ws.Range('E20').Value = '  +5.13%  '

$ws.Range('E21').Value = '  +4.11%  '

$ws.Range('D22').Value = "'375.03"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.18%  '

$ws.Range('D23').Value = "'0.573"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.55%  '

$ws.Range('D24').Value = '3.498.23'
$ws.Range('E24').Value = '  +2.10%  '

$ws.Range('E25').Value = '  +0.08%  '

$ws.Range('D26').Value = "'70.82"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.55%  '

$ws.Range('E27').Value = '  +12.95%  '

$ws.Range('D28').Value = "'1.64"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +17.94%  '

$ws.Range('D29').Value = "'7.73"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.59%  '

$ws.Range('E30').Value = '  -1.18%  '

$ws.Range('D31').Value = "'8.08"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.90%  '

$ws.Range('D32').Value = "'0.155"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.95%  '

$ws.Range('E33').Value = '  +3.55%  '

$ws.Range('E34').Value = '  -0.02%  '

$ws.Range('D35').Value = '3.397.90'
$ws.Range('E35').Value = '  +2.69%  '

$ws.Range('D36').Value = "'23.43"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.04%  '

$ws.Range('D37').Value = "'5.55"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.82%  '

$ws.Range('D38').Value = "'6.92"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.09%  '

$ws.Range('E39').Value = '  +6.16%  '

$ws.Range('D40').Value = "'164.22"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.07%  '

$ws.Range('D41').Value = "'0.0792"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.64%  '

$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.07%  '

$ws.Range('D43').Value = "'41.52"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.54%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = "'4.39"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.72%  '

$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').Value = "'1.20"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.47%  '

$ws.Range('E46').Value = '  +2.03%  '

$ws.Range('D47').Value = "'1.61"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.22%  '

$ws.Range('D48').Value = "'22.93"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.70%  '

$ws.Range('D49').Value = "'6.95"
$ws.Range('D49').Style = 'Normal'

$ws.Range('E50').Value = '  +14.63%  '

$ws.Range('D51').Value = "'2.42"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +16.43%  '
